$p = $ppt.ActivePresentation

# --- 1. Reorder slides: swap slide 2 and slide 3 -------------------------
# (drag slide 3 up to position 2, which pushes the old slide 2 down to 3)
$p.Slides.Item(3).MoveTo(2)

# --- 2. Refresh the "update automatically" date footer ---------------------
# The date placeholder (type 16 = ppPlaceholderDate) on the slide master and
# on every slide layout shows the day the deck was last touched; bump it from
# 7/26/2023 to 7/27/2023 everywhere it appears.
$master = $p.SlideMaster

foreach ($ph in $master.Shapes.Placeholders) {
    if ($ph.PlaceholderFormat.Type -eq 16) {
        $ph.TextFrame.TextRange.Text = "7/27/2023"
    }
}

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    foreach ($ph in $layout.Shapes.Placeholders) {
        if ($ph.PlaceholderFormat.Type -eq 16) {
            $ph.TextFrame.TextRange.Text = "7/27/2023"
        }
    }
}

Write-Output "Reordered slides 2/3 and refreshed $($master.CustomLayouts.Count + 1) date placeholders"
